$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Shared-string text edits (Volume/Number + report week dates) ----
$ws.Range("A8").Characters(21,2).Text = "48"
$ws.Range("C9").Characters(27,10).Text = "11/24/2025"
$ws.Range("C9").Characters(48,10).Text = "11/30/2025"

# ---- Numeric cell updates (rows 15-33) ----
# Row 15
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("M15").Value = -33.333333333333

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 400
$ws.Range("I16").Value = 47
$ws.Range("K16").Value = -11.320754716981
$ws.Range("M16").Value = -56.481481481481
$ws.Range("N16").Value = -83.392226148409

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 236
$ws.Range("J17").Value = 153
$ws.Range("K17").Value = 54.248366013071
$ws.Range("L17").Value = 51.282051282051
$ws.Range("M17").Value = 88.8
$ws.Range("N17").Value = -17.482517482517

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -37.5
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = -11.702127659574
$ws.Range("L18").Value = -7.777777777777
$ws.Range("M18").Value = -56.770833333333
$ws.Range("N18").Value = -93.386454183266

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 46.428571428571
$ws.Range("I19").Value = 382
$ws.Range("J19").Value = 434
$ws.Range("K19").Value = -11.981566820276
$ws.Range("L19").Value = -13.574660633484
$ws.Range("M19").Value = 0.526315789473
$ws.Range("N19").Value = -50.325097529258

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 166.666666666667
$ws.Range("I20").Value = 53
$ws.Range("K20").Value = -26.388888888888
$ws.Range("L20").Value = -47.524752475247
$ws.Range("M20").Value = -51.376146788990
$ws.Range("N20").Value = -97.938545313107

# Row 21
$ws.Range("C21").Value = 9
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -35.714285714285
$ws.Range("F21").Value = 76
$ws.Range("G21").Value = 58
$ws.Range("H21").Value = 31.034482758620
$ws.Range("I21").Value = 813
$ws.Range("J21").Value = 823
$ws.Range("K21").Value = -1.215066828675
$ws.Range("L21").Value = -3.900709219858
$ws.Range("M21").Value = -12.955032119914
$ws.Range("N21").Value = -84.320154291224

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -40
$ws.Range("I23").Value = 56
$ws.Range("J23").Value = 27
$ws.Range("K23").Value = 107.407407407407
$ws.Range("L23").Value = 75
$ws.Range("M23").Value = 154.545454545455

# Row 24
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 30
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 928
$ws.Range("J24").Value = 917
$ws.Range("K24").Value = 1.199563794983
$ws.Range("L24").Value = -6.920762286860
$ws.Range("M24").Value = -40.702875399361

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = -23.529411764705
$ws.Range("I25").Value = 474
$ws.Range("J25").Value = 449
$ws.Range("K25").Value = 5.567928730512
$ws.Range("L25").Value = 16.748768472906

# Row 26
$ws.Range("C26").Value = 5
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 19
$ws.Range("I26").Value = 331
$ws.Range("J26").Value = 332
$ws.Range("K26").Value = -0.301204819277
$ws.Range("L26").Value = 9.602649006622
$ws.Range("M26").Value = -33.8

# Row 27
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 0

# Row 28
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 200

# Row 33
$ws.Range("G33").Value = 1

# ---- Cells converting from numeric to special text markers ("0" / "***.*") ----
# C18 -> text "0" (style should match A18, s=13)
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("A18").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# D20 -> text "0" (style should match A20, s=13)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("A20").Copy()
$ws.Range("D20").PasteSpecial(-4122)

# E20 -> text "***.*" (style should match A20, s=13)
$ws.Range("E20").Value = "***.*"
$ws.Range("A20").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$excel.CutCopyMode = $false